# tag v 1.0.0 issue 2 issue 3
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the placeholder tokens from "{s-Name}" / "{s-Age}" style to
# the dotted "{s.Name}" / "{s.Age}" style used by the new template syntax.
$ws.Range("D7").Value = "{s.Name}"
$ws.Range("E7").Value = "{s.Age}"
